# Update header row text to the new (shortened) column headers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "MPA_Name"
$ws.Range("B1").Value = "Longitude"
$ws.Range("C1").Value = "Latitude"
$ws.Range("D1").Value = "Area"
$ws.Range("E1").Value = "State/UT"

# Match the selected cell recorded in the saved sheet view.
$ws.Range("G13").Select()
